$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns widths for L, M, N (12, 13, 14)
$ws.Columns.Item(12).ColumnWidth = 7.7109375
$ws.Columns.Item(13).ColumnWidth = 11.28515625
$ws.Columns.Item(14).ColumnWidth = 16.140625

# Header row for new "Grades" table (row 20)
$ws.Range("K20").Value = "Grades"
$ws.Range("L20").Value = "Factor"
$ws.Range("M20").Value = "% improved"
$ws.Range("N20").Value = "#Instructions"

# Data rows under the Grades table
$ws.Range("K21").Value = 100
$ws.Range("K22").Value = 90
$ws.Range("K23").Value = 80
$ws.Range("K24").Value = 70

# Bottom border under row 19 (closing off the original table) and around
# the new header row (row 20)
$ws.Range("F19:N19").Rows.RowHeight = 15.75
$ws.Range("K20:N20").Rows.RowHeight = 15.75

# Borders: outer box (medium) with thin separators between header cells
$headerRange = $ws.Range("K20:N20")
$headerRange.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$headerRange.Borders.Item(8).Weight = -4138  # xlMedium
$headerRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$headerRange.Borders.Item(9).Weight = -4138  # xlMedium
$headerRange.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$headerRange.Borders.Item(7).Weight = -4138  # xlMedium
$headerRange.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$headerRange.Borders.Item(10).Weight = -4138 # xlMedium
$headerRange.Borders.Item(11).LineStyle = 1  # xlInsideVertical
$headerRange.Borders.Item(11).Weight = 2     # xlThin

# Move the view so A3 is the top-left cell, and select L19 (matches target)
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("L19").Select()
